$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.4912210536349
$ws.Range("C2").Value = 7.079516359081339
$ws.Range("E2").Value = 16.26169844247173
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.689719431862343
$ws.Range("K2").Value = 13.84453735182456
$ws.Range("N2").Value = 21.62270149231703
$ws.Range("B3").Value = 14.20137136207102
$ws.Range("C3").Value = 6.82392300194945
$ws.Range("E3").Value = 15.35250948142855
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.693259905708472
$ws.Range("K3").Value = 13.64107894097801
$ws.Range("N3").Value = 21.64174064936995
$ws.Range("B4").Value = 14.02515504888891
$ws.Range("C4").Value = 6.664606629571034
$ws.Range("E4").Value = 14.77147776301505
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.695541450298191
$ws.Range("K4").Value = 13.51896853510846
$ws.Range("N4").Value = 21.65516783943433
$ws.Range("B5").Value = 13.95389361324201
$ws.Range("C5").Value = 6.599199982845853
$ws.Range("E5").Value = 14.5292495782045
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.696498390004387
$ws.Range("K5").Value = 13.46997782787656
$ws.Range("N5").Value = 21.66107488040349
$ws.Range("B6").Value = 13.94209696140963
$ws.Range("C6").Value = 6.588313570864923
$ws.Range("E6").Value = 14.48870754695591
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.696658934796372
$ws.Range("K6").Value = 13.46189128165816
$ws.Range("N6").Value = 21.66208199155207
$ws.Range("B7").Value = 14.02419163059901
$ws.Range("C7").Value = 6.663726333418224
$ws.Range("E7").Value = 14.76823266359479
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.695554245674172
$ws.Range("K7").Value = 13.51830462891288
$ws.Range("N7").Value = 21.65524574318691
$ws.Range("B8").Value = 14.39098135988876
$ws.Range("C8").Value = 6.99194785589312
$ws.Range("E8").Value = 15.95308972361363
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.690917907799909
$ws.Range("K8").Value = 13.77383893393437
$ws.Range("N8").Value = 21.62890479520794
$ws.Range("B9").Value = 15.11949922488199
$ws.Range("C9").Value = 7.612146351663402
$ws.Range("E9").Value = 18.09284780727973
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.682675191085282
$ws.Range("K9").Value = 14.29445235465257
$ws.Range("N9").Value = 21.59109871967752
$ws.Range("B10").Value = 15.65429785228634
$ws.Range("C10").Value = 8.048115177467672
$ws.Range("E10").Value = 19.71240002474015
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.677129530861658
$ws.Range("K10").Value = 14.68503357643522
$ws.Range("N10").Value = 21.57185600491345
$ws.Range("B11").Value = 15.89627692052382
$ws.Range("C11").Value = 8.241224211208987
$ws.Range("E11").Value = 20.4082400822486
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.674715878970088
$ws.Range("K11").Value = 14.86367167520345
$ws.Range("N11").Value = 21.56497437265686
$ws.Range("B12").Value = 15.98762522923459
$ws.Range("C12").Value = 8.313532704111482
$ws.Range("E12").Value = 20.66589358788439
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.673817457064806
$ws.Range("K12").Value = 14.93138976450192
$ws.Range("N12").Value = 21.56263921538601
$ws.Range("B13").Value = 15.96796584376084
$ws.Range("C13").Value = 8.297997134408986
$ws.Range("E13").Value = 20.61066273241322
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.674010257351394
$ws.Range("K13").Value = 14.91680331172749
$ws.Range("N13").Value = 21.56313006577616
$ws.Range("B14").Value = 15.90379842482883
$ws.Range("C14").Value = 8.247189820904259
$ws.Range("E14").Value = 20.42955435035663
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.674641653742612
$ws.Range("K14").Value = 14.86924182015463
$ws.Range("N14").Value = 21.56477682228405
$ws.Range("B15").Value = 15.86445428022691
$ws.Range("C15").Value = 8.215960516906041
$ws.Range("E15").Value = 20.31785996870638
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.675030427344273
$ws.Range("K15").Value = 14.84011643087484
$ws.Range("N15").Value = 21.56582081675903
$ws.Range("B16").Value = 15.6384493565478
$ws.Range("C16").Value = 8.035384117888945
$ws.Range("E16").Value = 19.66610523485372
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.677289454393995
$ws.Range("K16").Value = 14.67337277194686
$ws.Range("N16").Value = 21.57234354002378
$ws.Range("B17").Value = 15.49939853077789
$ws.Range("C17").Value = 7.923219666539645
$ws.Range("E17").Value = 19.25582418507218
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.678703156694513
$ws.Range("K17").Value = 14.57127922405659
$ws.Range("N17").Value = 21.57682560526221
$ws.Range("B18").Value = 15.41930293789333
$ws.Range("C18").Value = 7.858218237274427
$ws.Range("E18").Value = 19.01598932819554
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.679526556283018
$ws.Range("K18").Value = 14.51265202941681
$ws.Range("N18").Value = 21.57957967625559
$ws.Range("B19").Value = 15.39216683479216
$ws.Range("C19").Value = 7.836128361074109
$ws.Range("E19").Value = 18.93412268391262
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.679807113469031
$ws.Range("K19").Value = 14.49282011832209
$ws.Range("N19").Value = 21.58054235527004
$ws.Range("B20").Value = 15.5142135773167
$ws.Range("C20").Value = 7.935210725689346
$ws.Range("E20").Value = 19.29989771644986
$ws.Range("F20").Value = 20.2495528364879
$ws.Range("G20").Value = 3.678551603018756
$ws.Range("K20").Value = 14.58213798201347
$ws.Range("N20").Value = 21.57633024402872
$ws.Range("B21").Value = 15.92265439934251
$ws.Range("C21").Value = 8.262135851364217
$ws.Range("E21").Value = 20.48290864808347
$ws.Range("F21").Value = 21.46857628470567
$ws.Range("G21").Value = 3.674455775492151
$ws.Range("K21").Value = 14.88321035751948
$ws.Range("N21").Value = 21.56428576889477
$ws.Range("B22").Value = 16.18789824984568
$ws.Range("C22").Value = 8.471008425442209
$ws.Range("E22").Value = 21.22201603353881
$ws.Range("F22").Value = 22.22866616901555
$ws.Range("G22").Value = 3.67186965033526
$ws.Range("K22").Value = 15.08036708830419
$ws.Range("N22").Value = 21.55799297205185
$ws.Range("B23").Value = 16.04651828446485
$ws.Range("C23").Value = 8.359987915213924
$ws.Range("E23").Value = 20.83064513180085
$ws.Range("F23").Value = 21.82633154475864
$ws.Range("G23").Value = 3.673241648846009
$ws.Range("K23").Value = 14.97512675840279
$ws.Range("N23").Value = 21.56120656221926
$ws.Range("B24").Value = 15.50751616527185
$ws.Range("C24").Value = 7.929791171687558
$ws.Range("E24").Value = 19.27998440644598
$ws.Range("F24").Value = 20.22900810905294
$ws.Range("G24").Value = 3.678620087275534
$ws.Range("K24").Value = 14.57722852146893
$ws.Range("N24").Value = 21.57655364479955
$ws.Range("B25").Value = 14.92206605213655
$ws.Range("C25").Value = 7.447471229087428
$ws.Range("E25").Value = 17.5306229213648
$ws.Range("F25").Value = 18.34778573295697
$ws.Range("G25").Value = 3.684814923419621
$ws.Range("K25").Value = 14.15191702507185
$ws.Range("N25").Value = 21.59983424328531
